$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "D=<n>" header labels to "#features=<n>" (B1:G1).
$ws.Range("B1").Value = "#features=50"
$ws.Range("C1").Value = "#features=100"
$ws.Range("D1").Value = "#features=150"
$ws.Range("E1").Value = "#features=200"
$ws.Range("F1").Value = "#features=250"
$ws.Range("G1").Value = "#features=315"

# Add normalized throughput: divide the throughput table (B2:G7) by 1000.
$dataRange = $ws.Range("B2:G7")
foreach ($cell in $dataRange.Cells) {
    $val = $cell.Value()
    $cell.Value = $val / 1000
}

# Clean up the stray formatted-but-empty cells below the table.
$ws.Range("B11:G16").Clear()

# Update the selection left on the sheet.
$ws.Range("C11").Select()
